# Prepocessing Player Stats for year 2019-20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Cleaning files" rows (F2:F4) as Completed
$ws.Range("F2").Value = "Completed"
$ws.Range("F3").Value = "Completed"
$ws.Range("F4").Value = "Completed"

# Update the selection / scroll position on Sheet1
$ws.Range("F1").Select()
$excel.ActiveWindow.ScrollRow = 2

# Update workbook window size/position
$excel.Width = 20760
$excel.Height = 18680
$excel.Left = 6160
$excel.Top = 1980
